$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 12: rename model combo, add submission score ---
$ws.Cells.Item(12, 2).Value = "model_floodwater_unet_pc_augm_diceloss + 2"
$ws.Cells.Item(12, 4).Value = 0.908

# --- Update row 14: was the diceloss_2 model row, becomes diceloss_3 with new data ---
$ws.Cells.Item(14, 2).Value = "model_floodwater_unet_pc_augm_diceloss_3"
$ws.Cells.Item(14, 3).Value = 0.681
$ws.Cells.Item(14, 5).Value = "['hbe', 'jja']"

# --- Add new row 15 ---
$ws.Cells.Item(15, 1).Value = 44457
$ws.Cells.Item(15, 1).NumberFormat = "m/d/yyyy"
$ws.Cells.Item(15, 2).Value = "model_floodwater_unet_pc_augm_diceloss 1 +2 + 3"
$ws.Cells.Item(15, 3).Value = 0.718

# --- Add new row 16 ---
$ws.Cells.Item(16, 1).Value = 44457
$ws.Cells.Item(16, 1).NumberFormat = "m/d/yyyy"
$ws.Cells.Item(16, 2).Value = "model_floodwater_unet_pc_augm_diceloss_4"
$ws.Cells.Item(16, 3).Value = 0.699
$ws.Cells.Item(16, 5).Value = "['pxs', 'tnp'] "

# --- Add new row 17 ---
$ws.Cells.Item(17, 1).Value = 44457
$ws.Cells.Item(17, 1).NumberFormat = "m/d/yyyy"
$ws.Cells.Item(17, 2).Value = "model_floodwater_unet_pc_augm_diceloss 1 +2 + 3 + 4"
$ws.Cells.Item(17, 3).Value = 0.717

# --- Add new row 18 ---
$ws.Cells.Item(18, 1).Value = 44459
$ws.Cells.Item(18, 1).NumberFormat = "m/d/yyyy"
$ws.Cells.Item(18, 2).Value = "model_floodwater_unet_pc_augm_diceloss 1 +2 + 4"

# --- Update sheet view: scroll + selection to match final saved state ---
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("C18").Select()
